# edit.ps1
# Commit: "Created new caesar cypher" - extends the ETA-tracking data table
# (rows 21-48) on Sheet1, updates the combo chart's data ranges to match,
# repositions/resizes the chart, and refreshes the sheet view (zoom +
# selection) to reflect where the user was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Extend the data table: rows 21-48 continue the same formula pattern
#    already used in rows 3-20 (A: running index, B: raw elapsed-word
#    count input, C: running total, E: marginal delta, F: constant 5000,
#    G: ETA-hours formula). B is the only literal input; everything else
#    is a formula so Excel (and this engine) compute matching values.
# ---------------------------------------------------------------------
$ws.Range("A21").Formula = "=A20+1"
$ws.Range("B21").Value = 6102
$ws.Range("C21").Formula = "=C20+5000"
$ws.Range("E21").Formula = "=B21-B20"
$ws.Range("F21").Formula = "=5000"
$ws.Range("G21").Formula = "=B21*(235886/C21-1)/3600"
$ws.Range("A22").Formula = "=A21+1"
$ws.Range("B22").Value = 6285
$ws.Range("C22").Formula = "=C21+5000"
$ws.Range("E22").Formula = "=B22-B21"
$ws.Range("F22").Formula = "=5000"
$ws.Range("G22").Formula = "=B22*(235886/C22-1)/3600"
$ws.Range("A23").Formula = "=A22+1"
$ws.Range("B23").Value = 6420
$ws.Range("C23").Formula = "=C22+5000"
$ws.Range("E23").Formula = "=B23-B22"
$ws.Range("F23").Formula = "=5000"
$ws.Range("G23").Formula = "=B23*(235886/C23-1)/3600"
$ws.Range("A24").Formula = "=A23+1"
$ws.Range("B24").Value = 6657
$ws.Range("C24").Formula = "=C23+5000"
$ws.Range("E24").Formula = "=B24-B23"
$ws.Range("F24").Formula = "=5000"
$ws.Range("G24").Formula = "=B24*(235886/C24-1)/3600"
$ws.Range("A25").Formula = "=A24+1"
$ws.Range("B25").Value = 6965
$ws.Range("C25").Formula = "=C24+5000"
$ws.Range("E25").Formula = "=B25-B24"
$ws.Range("F25").Formula = "=5000"
$ws.Range("G25").Formula = "=B25*(235886/C25-1)/3600"
$ws.Range("A26").Formula = "=A25+1"
$ws.Range("B26").Value = 7406
$ws.Range("C26").Formula = "=C25+5000"
$ws.Range("E26").Formula = "=B26-B25"
$ws.Range("F26").Formula = "=5000"
$ws.Range("G26").Formula = "=B26*(235886/C26-1)/3600"
$ws.Range("A27").Formula = "=A26+1"
$ws.Range("B27").Value = 7712
$ws.Range("C27").Formula = "=C26+5000"
$ws.Range("E27").Formula = "=B27-B26"
$ws.Range("F27").Formula = "=5000"
$ws.Range("G27").Formula = "=B27*(235886/C27-1)/3600"
$ws.Range("A28").Formula = "=A27+1"
$ws.Range("B28").Value = 7903
$ws.Range("C28").Formula = "=C27+5000"
$ws.Range("E28").Formula = "=B28-B27"
$ws.Range("F28").Formula = "=5000"
$ws.Range("G28").Formula = "=B28*(235886/C28-1)/3600"
$ws.Range("A29").Formula = "=A28+1"
$ws.Range("B29").Value = 8063
$ws.Range("C29").Formula = "=C28+5000"
$ws.Range("E29").Formula = "=B29-B28"
$ws.Range("F29").Formula = "=5000"
$ws.Range("G29").Formula = "=B29*(235886/C29-1)/3600"
$ws.Range("A30").Formula = "=A29+1"
$ws.Range("B30").Value = 8307
$ws.Range("C30").Formula = "=C29+5000"
$ws.Range("E30").Formula = "=B30-B29"
$ws.Range("F30").Formula = "=5000"
$ws.Range("G30").Formula = "=B30*(235886/C30-1)/3600"
$ws.Range("A31").Formula = "=A30+1"
$ws.Range("B31").Value = 8636
$ws.Range("C31").Formula = "=C30+5000"
$ws.Range("E31").Formula = "=B31-B30"
$ws.Range("F31").Formula = "=5000"
$ws.Range("G31").Formula = "=B31*(235886/C31-1)/3600"
$ws.Range("A32").Formula = "=A31+1"
$ws.Range("B32").Value = 9060
$ws.Range("C32").Formula = "=C31+5000"
$ws.Range("E32").Formula = "=B32-B31"
$ws.Range("F32").Formula = "=5000"
$ws.Range("G32").Formula = "=B32*(235886/C32-1)/3600"
$ws.Range("H32").Value = "8:01PM"
$ws.Range("A33").Formula = "=A32+1"
$ws.Range("B33").Value = 9594
$ws.Range("C33").Formula = "=C32+5000"
$ws.Range("E33").Formula = "=B33-B32"
$ws.Range("F33").Formula = "=5000"
$ws.Range("G33").Formula = "=B33*(235886/C33-1)/3600"
$ws.Range("A34").Formula = "=A33+1"
$ws.Range("B34").Value = 10225
$ws.Range("C34").Formula = "=C33+5000"
$ws.Range("E34").Formula = "=B34-B33"
$ws.Range("F34").Formula = "=5000"
$ws.Range("G34").Formula = "=B34*(235886/C34-1)/3600"
$ws.Range("A35").Formula = "=A34+1"
$ws.Range("B35").Value = 10464
$ws.Range("C35").Formula = "=C34+5000"
$ws.Range("E35").Formula = "=B35-B34"
$ws.Range("F35").Formula = "=5000"
$ws.Range("G35").Formula = "=B35*(235886/C35-1)/3600"
$ws.Range("A36").Formula = "=A35+1"
$ws.Range("B36").Value = 10682
$ws.Range("C36").Formula = "=C35+5000"
$ws.Range("E36").Formula = "=B36-B35"
$ws.Range("F36").Formula = "=5000"
$ws.Range("G36").Formula = "=B36*(235886/C36-1)/3600"
$ws.Range("A37").Formula = "=A36+1"
$ws.Range("B37").Value = 11004
$ws.Range("C37").Formula = "=C36+5000"
$ws.Range("E37").Formula = "=B37-B36"
$ws.Range("F37").Formula = "=5000"
$ws.Range("G37").Formula = "=B37*(235886/C37-1)/3600"
$ws.Range("A38").Formula = "=A37+1"
$ws.Range("B38").Value = 11393
$ws.Range("C38").Formula = "=C37+5000"
$ws.Range("E38").Formula = "=B38-B37"
$ws.Range("F38").Formula = "=5000"
$ws.Range("G38").Formula = "=B38*(235886/C38-1)/3600"
$ws.Range("A39").Formula = "=A38+1"
$ws.Range("B39").Value = 11879
$ws.Range("C39").Formula = "=C38+5000"
$ws.Range("E39").Formula = "=B39-B38"
$ws.Range("F39").Formula = "=5000"
$ws.Range("G39").Formula = "=B39*(235886/C39-1)/3600"
$ws.Range("A40").Formula = "=A39+1"
$ws.Range("B40").Value = 12518
$ws.Range("C40").Formula = "=C39+5000"
$ws.Range("E40").Formula = "=B40-B39"
$ws.Range("F40").Formula = "=5000"
$ws.Range("G40").Formula = "=B40*(235886/C40-1)/3600"
$ws.Range("H40").Value = "8:59PM"
$ws.Range("A41").Formula = "=A40+1"
$ws.Range("B41").Value = 13275
$ws.Range("C41").Formula = "=C40+5000"
$ws.Range("E41").Formula = "=B41-B40"
$ws.Range("F41").Formula = "=5000"
$ws.Range("G41").Formula = "=B41*(235886/C41-1)/3600"
$ws.Range("A42").Formula = "=A41+1"
$ws.Range("B42").Value = 13724
$ws.Range("C42").Formula = "=C41+5000"
$ws.Range("E42").Formula = "=B42-B41"
$ws.Range("F42").Formula = "=5000"
$ws.Range("G42").Formula = "=B42*(235886/C42-1)/3600"
$ws.Range("A43").Formula = "=A42+1"
$ws.Range("B43").Value = 13994
$ws.Range("C43").Formula = "=C42+5000"
$ws.Range("E43").Formula = "=B43-B42"
$ws.Range("F43").Formula = "=5000"
$ws.Range("G43").Formula = "=B43*(235886/C43-1)/3600"
$ws.Range("A44").Formula = "=A43+1"
$ws.Range("B44").Value = 14383
$ws.Range("C44").Formula = "=C43+5000"
$ws.Range("E44").Formula = "=B44-B43"
$ws.Range("F44").Formula = "=5000"
$ws.Range("G44").Formula = "=B44*(235886/C44-1)/3600"
$ws.Range("A45").Formula = "=A44+1"
$ws.Range("B45").Value = 14538
$ws.Range("C45").Formula = "=C44+5000"
$ws.Range("E45").Formula = "=B45-B44"
$ws.Range("F45").Formula = "=5000"
$ws.Range("G45").Formula = "=B45*(235886/C45-1)/3600"
$ws.Range("A46").Formula = "=A45+1"
$ws.Range("B46").Value = 14740
$ws.Range("C46").Formula = "=C45+5000"
$ws.Range("E46").Formula = "=B46-B45"
$ws.Range("F46").Formula = "=5000"
$ws.Range("G46").Formula = "=B46*(235886/C46-1)/3600"
$ws.Range("A47").Formula = "=A46+1"
$ws.Range("B47").Value = 15021
$ws.Range("C47").Formula = "=C46+5000"
$ws.Range("E47").Formula = "=B47-B46"
$ws.Range("F47").Formula = "=5000"
$ws.Range("G47").Formula = "=B47*(235886/C47-1)/3600"
$ws.Range("A48").Formula = "=A47+1"
$ws.Range("B48").Value = 15190
$ws.Range("C48").Formula = "=C47+5000"
$ws.Range("E48").Formula = "=B48-B47"
$ws.Range("F48").Formula = "=5000"
$ws.Range("G48").Formula = "=B48*(235886/C48-1)/3600"

# ---------------------------------------------------------------------
# 2. Re-point the chart's two series at the now-extended ranges
#    (E3:E20 -> E3:E48, G3:G48 -> G3:G48) so the bar/line combo chart
#    picks up the new rows.
# ---------------------------------------------------------------------
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$chart.SeriesCollection(1).Formula = '=SERIES("Marginal Time Elapsed",,Sheet1!$E$3:$E$48,1)'
$chart.SeriesCollection(2).Formula = '=SERIES("Estimated time remaining",,Sheet1!$G$3:$G$48,2)'

# ---------------------------------------------------------------------
# 3. Move/resize the chart to its new anchor position (shifted right and
#    down slightly, and a touch wider) to make room next to the grown
#    table.
# ---------------------------------------------------------------------
$co.Left = 704.4421576649
$co.Top = 56.8455118110
$co.Width = 433.2389566929
$co.Height = 216.0

# ---------------------------------------------------------------------
# 4. Sheet view: zoomed out slightly and selection left on the last cell
#    that was filled in.
# ---------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 85
$ws.Range("B48").Select()
